$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 21.02666664123535
$ws.Range("E2").Value = 23.70333290100098
$ws.Range("F2").Value = 24.3799991607666
$ws.Range("G2").Value = 20.21833229064941
$ws.Range("H2").Value = 676870933
$ws.Range("I2").Value = "PANW"

$ws.Range("D3").Value = 25.01166725158692
$ws.Range("E3").Value = 28.24833297729492
$ws.Range("F3").Value = 28.30666732788086
$ws.Range("G3").Value = 24.06999969482422
$ws.Range("H3").Value = 676870933
$ws.Range("I3").Value = "PANW"

$ws.Range("D4").Value = 31.0583324432373
$ws.Range("E4").Value = 27.3700008392334
$ws.Range("F4").Value = 31.83333206176757
$ws.Range("G4").Value = 23.3983325958252
$ws.Range("H4").Value = 676870933
$ws.Range("I4").Value = "PANW"

$ws.Range("D5").Value = 26.96666717529297
$ws.Range("E5").Value = 31.22333335876465
$ws.Range("F5").Value = 31.38333320617676
$ws.Range("G5").Value = 25.47999954223633
$ws.Range("H5").Value = 676870933
$ws.Range("I5").Value = "PANW"

$ws.Range("D6").Value = 24.64666748046875
$ws.Range("E6").Value = 24.13166618347168
$ws.Range("F6").Value = 25.65500068664551
$ws.Range("G6").Value = 18.51499938964844
$ws.Range("H6").Value = 676870933
$ws.Range("I6").Value = "PANW"

$ws.Range("D7").Value = 25.27499961853028
$ws.Range("E7").Value = 21.74333381652832
$ws.Range("F7").Value = 25.33166694641113
$ws.Range("G7").Value = 21.375
$ws.Range("H7").Value = 676870933
$ws.Range("I7").Value = "PANW"

$ws.Range("D8").Value = 21.80999946594238
$ws.Range("E8").Value = 22.19499969482422
$ws.Range("F8").Value = 24.04166793823243
$ws.Range("G8").Value = 20.79000091552734
$ws.Range("H8").Value = 676870933
$ws.Range("I8").Value = "PANW"

$ws.Range("D9").Value = 25.625
$ws.Range("E9").Value = 22.39500045776367
$ws.Range("F9").Value = 27.61499977111816
$ws.Range("G9").Value = 22.02000045776367
$ws.Range("H9").Value = 676870933
$ws.Range("I9").Value = "PANW"

$ws.Range("D10").Value = 24.66666793823243
$ws.Range("E10").Value = 25.31666755676269
$ws.Range("F10").Value = 26.27499961853028
$ws.Range("G10").Value = 24.1016674041748
$ws.Range("H10").Value = 676870933
$ws.Range("I10").Value = "PANW"

$ws.Range("D11").Value = 18.125
$ws.Range("E11").Value = 19.76499938964844
$ws.Range("F11").Value = 20.34666633605957
$ws.Range("G11").Value = 18.02499961853028
$ws.Range("H11").Value = 676870933
$ws.Range("I11").Value = "PANW"

$ws.Range("D12").Value = 22.09666633605957
$ws.Range("E12").Value = 22.11499977111816
$ws.Range("F12").Value = 22.47999954223633
$ws.Range("G12").Value = 21.09333229064941
$ws.Range("H12").Value = 676870933
$ws.Range("I12").Value = "PANW"

$ws.Range("D13").Value = 24.59666633605957
$ws.Range("E13").Value = 24.29166793823243
$ws.Range("F13").Value = 26.14166641235352
$ws.Range("G13").Value = 22.64166641235352
$ws.Range("H13").Value = 676870933
$ws.Range("I13").Value = "PANW"

$ws.Range("D14").Value = 26.10000038146973
$ws.Range("E14").Value = 28.89500045776367
$ws.Range("F14").Value = 29.85833358764648
$ws.Range("G14").Value = 24.73500061035156
$ws.Range("H14").Value = 676870933
$ws.Range("I14").Value = "PANW"

$ws.Range("D15").Value = 31.95833206176757
$ws.Range("E15").Value = 34.68166732788086
$ws.Range("F15").Value = 35.28499984741211
$ws.Range("G15").Value = 31.63999938964844
$ws.Range("H15").Value = 676870933
$ws.Range("I15").Value = "PANW"

$ws.Range("D16").Value = 33.4283332824707
$ws.Range("E16").Value = 38.52500152587891
$ws.Range("F16").Value = 38.71833419799805
$ws.Range("G16").Value = 33.00333404541016
$ws.Range("H16").Value = 676870933
$ws.Range("I16").Value = "PANW"

$ws.Range("D17").Value = 30.75
$ws.Range("E17").Value = 28.82500076293945
$ws.Range("F17").Value = 32.04499816894531
$ws.Range("G17").Value = 26.68000030517578
$ws.Range("H17").Value = 676870933
$ws.Range("I17").Value = "PANW"

$ws.Range("D18").Value = 35.84999847412109
$ws.Range("E18").Value = 41.04499816894531
$ws.Range("F18").Value = 43.43833160400391
$ws.Range("G18").Value = 35.73333358764648
$ws.Range("H18").Value = 676870933
$ws.Range("I18").Value = "PANW"

$ws.Range("D19").Value = 41.58166885375977
$ws.Range("E19").Value = 33.35666656494141
$ws.Range("F19").Value = 41.82333374023438
$ws.Range("G19").Value = 32.8466682434082
$ws.Range("H19").Value = 676870933
$ws.Range("I19").Value = "PANW"

$ws.Range("D20").Value = 37.77500152587891
$ws.Range("E20").Value = 33.93666839599609
$ws.Range("F20").Value = 38.32333374023438
$ws.Range("G20").Value = 32.0283317565918
$ws.Range("H20").Value = 676870933
$ws.Range("I20").Value = "PANW"

$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 37.86999893188477
$ws.Range("F21").Value = 41.82500076293945
$ws.Range("G21").Value = 35.81666564941406
$ws.Range("H21").Value = 676870933
$ws.Range("I21").Value = "PANW"

$ws.Range("D22").Value = 39.14333343505859
$ws.Range("E22").Value = 30.77000045776367
$ws.Range("F22").Value = 41.85166549682617
$ws.Range("G22").Value = 30.17666625976562
$ws.Range("H22").Value = 676870933
$ws.Range("I22").Value = "PANW"

$ws.Range("D23").Value = 32.16666793823242
$ws.Range("E23").Value = 39.21166610717773
$ws.Range("F23").Value = 40.52333450317383
$ws.Range("G23").Value = 31.75833320617676
$ws.Range("H23").Value = 676870933
$ws.Range("I23").Value = "PANW"

$ws.Range("D24").Value = 42.6533317565918
$ws.Range("E24").Value = 42.90166854858398
$ws.Range("F24").Value = 45.83833312988281
$ws.Range("G24").Value = 41.78666687011719
$ws.Range("H24").Value = 676870933
$ws.Range("I24").Value = "PANW"

$ws.Range("D25").Value = 37.13166809082031
$ws.Range("E25").Value = 48.98666763305664
$ws.Range("F25").Value = 49.75
$ws.Range("G25").Value = 36.82833480834961
$ws.Range("H25").Value = 676870933
$ws.Range("I25").Value = "PANW"

$ws.Range("D26").Value = 59.14833450317383
$ws.Range("E26").Value = 59.71833419799805
$ws.Range("F26").Value = 67.16666412353516
$ws.Range("G26").Value = 58.66666793823242
$ws.Range("H26").Value = 676870933
$ws.Range("I26").Value = "PANW"

$ws.Range("D27").Value = 59.26166534423828
$ws.Range("E27").Value = 60.54166793823242
$ws.Range("F27").Value = 62.38166809082031
$ws.Range("G27").Value = 53.70500183105469
$ws.Range("H27").Value = 676870933
$ws.Range("I27").Value = "PANW"

$ws.Range("D28").Value = 66.66666412353516
$ws.Range("E28").Value = 76.83999633789062
$ws.Range("F28").Value = 77.48332977294922
$ws.Range("G28").Value = 59.72833251953125
$ws.Range("H28").Value = 676870933
$ws.Range("I28").Value = "PANW"

$ws.Range("D29").Value = 85.08499908447266
$ws.Range("E29").Value = 91.15666961669922
$ws.Range("F29").Value = 93.25666809082033
$ws.Range("G29").Value = 80.87666320800781
$ws.Range("H29").Value = 676870933
$ws.Range("I29").Value = "PANW"

$ws.Range("D30").Value = 86.30999755859375
$ws.Range("E30").Value = 99.04166412353516
$ws.Range("F30").Value = 99.21666717529295
$ws.Range("G30").Value = 77.75
$ws.Range("H30").Value = 676870933
$ws.Range("I30").Value = "PANW"

$ws.Range("D31").Value = 93.0433349609375
$ws.Range("E31").Value = 83.79666900634766
$ws.Range("F31").Value = 95.47833251953124
$ws.Range("G31").Value = 70.25833129882812
$ws.Range("H31").Value = 676870933
$ws.Range("I31").Value = "PANW"

$ws.Range("D32").Value = 83.2933349609375
$ws.Range("E32").Value = 92.80166625976562
$ws.Range("F32").Value = 96.46499633789062
$ws.Range("G32").Value = 80.67832946777344
$ws.Range("H32").Value = 676870933
$ws.Range("I32").Value = "PANW"

$ws.Range("D33").Value = 86.97000122070312
$ws.Range("E33").Value = 84.94999694824219
$ws.Range("F33").Value = 88.05000305175781
$ws.Range("G33").Value = 70.03500366210938
$ws.Range("H33").Value = 676870933
$ws.Range("I33").Value = "PANW"

$ws.Range("D34").Value = 78.67500305175781
$ws.Range("E34").Value = 94.18499755859376
$ws.Range("F34").Value = 96.23500061035156
$ws.Range("G34").Value = 77.05500030517578
$ws.Range("H34").Value = 676870933
$ws.Range("I34").Value = "PANW"

$ws.Range("D35").Value = 90.33000183105467
$ws.Range("E35").Value = 106.6949996948242
$ws.Range("F35").Value = 109.75
$ws.Range("G35").Value = 88.15000152587891
$ws.Range("H35").Value = 676870933
$ws.Range("I35").Value = "PANW"

$ws.Range("D36").Value = 124.9449996948242
$ws.Range("E36").Value = 121.6500015258789
$ws.Range("F36").Value = 127.0149993896484
$ws.Range("G36").Value = 100.5849990844727
$ws.Range("H36").Value = 676870933
$ws.Range("I36").Value = "PANW"

$ws.Range("D37").Value = 122.5
$ws.Range("E37").Value = 147.5449981689453
$ws.Range("F37").Value = 148.0800018310547
$ws.Range("G37").Value = 116.9049987792969
$ws.Range("H37").Value = 676870933
$ws.Range("I37").Value = "PANW"

$ws.Range("D38").Value = 169.5
$ws.Range("E38").Value = 155.2749938964844
$ws.Range("F38").Value = 190.4199981689453
$ws.Range("G38").Value = 130.0449981689453
$ws.Range("H38").Value = 676870933
$ws.Range("I38").Value = "PANW"

$ws.Range("D39").Value = 144.4400024414062
$ws.Range("E39").Value = 147.4550018310547
$ws.Range("F39").Value = 162.3399963378906
$ws.Range("G39").Value = 142
$ws.Range("H39").Value = 676870933
$ws.Range("I39").Value = "PANW"

$ws.Range("D40").Value = 161.7050018310547
$ws.Range("E40").Value = 181.3600006103516
$ws.Range("F40").Value = 187.6849975585937
$ws.Range("G40").Value = 142.0099945068359
$ws.Range("H40").Value = 676870933
$ws.Range("I40").Value = "PANW"

$ws.Range("D41").Value = 179.5149993896484
$ws.Range("E41").Value = 193.9100036621093
$ws.Range("F41").Value = 204.2649993896484
$ws.Range("G41").Value = 178
$ws.Range("H41").Value = 676870933
$ws.Range("I41").Value = "PANW"

$ws.Range("D42").Value = 181.5599975585937
$ws.Range("E42").Value = 190.4299926757812
$ws.Range("F42").Value = 208.3899993896484
$ws.Range("G42").Value = 180.1199951171875
$ws.Range("H42").Value = 676870933
$ws.Range("I42").Value = "PANW"

$ws.Range("D43").Value = 188.6399993896484
$ws.Range("E43").Value = 192.4199981689453
$ws.Range("F43").Value = 195.4199981689453
$ws.Range("G43").Value = 178.6399993896484
$ws.Range("H43").Value = 676870933
$ws.Range("I43").Value = "PANW"

$ws.Range("D44").Value = 173.1300048828125
$ws.Range("E44").Value = 190.5200042724609
$ws.Range("F44").Value = 192.7899932861328
$ws.Range("G44").Value = 165.2100067138672
$ws.Range("H44").Value = 676870933
$ws.Range("I44").Value = "PANW"
